# Apply updated monthly indicator values ("Actualización desde MV -datos-")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D62").Value = 1.8
$ws.Range("C63").Value = 83.3
$ws.Range("D67").Value = 0.1
$ws.Range("C68").Value = 83.59999999999999
$ws.Range("C70").Value = 82.40000000000001
$ws.Range("D71").Value = -1.6
$ws.Range("D72").Value = -1.1
$ws.Range("C74").Value = 80.90000000000001
$ws.Range("C75").Value = 80.7
$ws.Range("D75").Value = -0.9
$ws.Range("C76").Value = 80.8
$ws.Range("C82").Value = 82
$ws.Range("D85").Value = 1.2
$ws.Range("D86").Value = 1.1
$ws.Range("D87").Value = 1.5
$ws.Range("C94").Value = 87.8
$ws.Range("D94").Value = 2.4
$ws.Range("D97").Value = 1.2
$ws.Range("D98").Value = 1.5
$ws.Range("D108").Value = 1.2
$ws.Range("C109").Value = 94.09999999999999
$ws.Range("C110").Value = 93.8
$ws.Range("D111").Value = 2.3
$ws.Range("D113").Value = 1.9
$ws.Range("D114").Value = 1.6
$ws.Range("C117").Value = 96.8
$ws.Range("C119").Value = 96.7
$ws.Range("D122").Value = 1.2
$ws.Range("D123").Value = 1.2
$ws.Range("D132").Value = 0.1
$ws.Range("D134").Value = 0.1
$ws.Range("D135").Value = -0.2
$ws.Range("C136").Value = 101.2
$ws.Range("D145").Value = 1.1
$ws.Range("D150").Value = 0.7
$ws.Range("D152").Value = 1
$ws.Range("C153").Value = 103.1
$ws.Range("D156").Value = 0.8
$ws.Range("D157").Value = 0.9
$ws.Range("D159").Value = 0.5
$ws.Range("C163").Value = 106
$ws.Range("C165").Value = 105.8
$ws.Range("D165").Value = 0.4
$ws.Range("C169").Value = 105.6
$ws.Range("D169").Value = 0.2
$ws.Range("C170").Value = 105.8
$ws.Range("D170").Value = 0
$ws.Range("D171").Value = -0.5
$ws.Range("C173").Value = 105.9
$ws.Range("D173").Value = -0.7
$ws.Range("C174").Value = 106.4
$ws.Range("C175").Value = 106.9
$ws.Range("C177").Value = 108.1
$ws.Range("D177").Value = 1.9
$ws.Range("C179").Value = 108.7
$ws.Range("D179").Value = 1.4
$ws.Range("C180").Value = 109.4
$ws.Range("D180").Value = 1.3
$ws.Range("D181").Value = 0.8
$ws.Range("C182").Value = 110
$ws.Range("C183").Value = 110.3
$ws.Range("D183").Value = 0.8
$ws.Range("D184").Value = 1.2
$ws.Range("C185").Value = 111.5
$ws.Range("D185").Value = 1.3
$ws.Range("C186").Value = 111.6
$ws.Range("D186").Value = 1.4
$ws.Range("C187").Value = 111.3
$ws.Range("D187").Value = 1
$ws.Range("C188").Value = 111.3
$ws.Range("D188").Value = 0.5
$ws.Range("D189").Value = -0.2
$ws.Range("D190").Value = -0.5
$ws.Range("C191").Value = 111.7
$ws.Range("D191").Value = -0.3
$ws.Range("C192").Value = 113
$ws.Range("D192").Value = 0.8
$ws.Range("C193").Value = 111.7
$ws.Range("D193").Value = 1
$ws.Range("C194").Value = 112.2
$ws.Range("D194").Value = 1.1
$ws.Range("C197").Value = 113.3
$ws.Range("D197").Value = 0
$ws.Range("C198").Value = 113.9
$ws.Range("D198").Value = 1.3
$ws.Range("C199").Value = 113.1
$ws.Range("D199").Value = 1.3
$ws.Range("C200").Value = 113.3
$ws.Range("D200").Value = 1
$ws.Range("C202").Value = 114.1
$ws.Range("D202").Value = 0.5
$ws.Range("C203").Value = 107.9
$ws.Range("D203").Value = -1.1
$ws.Range("C204").Value = 108.6
$ws.Range("D204").Value = -3
$ws.Range("C205").Value = 111.9
$ws.Range("D205").Value = -3.9
$ws.Range("C206").Value = 113.7
$ws.Range("D206").Value = -0.7
$ws.Range("C208").Value = 107.8
$ws.Range("D208").Value = 2.1
$ws.Range("C209").Value = 99.40000000000001
$ws.Range("D209").Value = -3.9
$ws.Range("C210").Value = 97.40000000000001
$ws.Range("D210").Value = -10.3
$ws.Range("C211").Value = 97
$ws.Range("D211").Value = -12.4
$ws.Range("C212").Value = 99.2
$ws.Range("C214").Value = 106.9
$ws.Range("D214").Value = 5
$ws.Range("C215").Value = 107.2
$ws.Range("D215").Value = 7.8
$ws.Range("C216").Value = 108.9
$ws.Range("D216").Value = 8.199999999999999
$ws.Range("C217").Value = 111.8
$ws.Range("D217").Value = 6.3
$ws.Range("C218").Value = 112.7
$ws.Range("C219").Value = 114.2
$ws.Range("D219").Value = 4.8
$ws.Range("C220").Value = 113.1
$ws.Range("D220").Value = 3.6
$ws.Range("C221").Value = 111.8
$ws.Range("D221").Value = 1.7
$ws.Range("C222").Value = 115.4
$ws.Range("D222").Value = 0.5
$ws.Range("C223").Value = 117.5
$ws.Range("D223").Value = 1.4
$ws.Range("C224").Value = 119.2
$ws.Range("D224").Value = 3.8

# New row 225: 01-08-2021 data point
# Force column A to stay text (avoid auto date-parsing) like the rest of the date column,
# then restore the default (Normal) style so no stray number format sticks to the cell.
$ws.Range("A225").NumberFormat = "@"
$ws.Range("A225").Value = "01-08-2021"
$ws.Range("A225").Style = "Normal"
$ws.Range("B225").Value = 120.3
$ws.Range("C225").Value = 120.4
$ws.Range("D225").Value = 4.9

